$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 29 - shifts rows 29:79 down to 30:80
$ws.Rows.Item(29).Insert()

# Populate the newly inserted row 29 with the new record
$ws.Cells.Item(29, 1).Value = 8
$ws.Cells.Item(29, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(29, 3).Value = "Coquimbo"
$ws.Cells.Item(29, 4).Value = 44883
$ws.Cells.Item(29, 5).Value = 4
$ws.Cells.Item(29, 6).Value = 100114007
$ws.Cells.Item(29, 7).Value = "Jengibre"
$ws.Cells.Item(29, 8).Value = "Sin especificar"
$ws.Cells.Item(29, 9).Value = "Primera"
$ws.Cells.Item(29, 10).Value = 480
$ws.Cells.Item(29, 11).Value = 13500
$ws.Cells.Item(29, 12).Value = 14000
$ws.Cells.Item(29, 13).Value = 13750
$ws.Cells.Item(29, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(29, 15).Value = "Perú"
$ws.Cells.Item(29, 16).Value = 1058
$ws.Cells.Item(29, 17).Value = 13
$ws.Cells.Item(29, 18).Value = "Hortaliza"
